$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 6, 9, 11, 13)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "12,5%"
    $ws.Range("F$r").Value = "87,5%"
}
